# Work on Input Capture Module
# Add a new measurement row (150 kHz) to the accuracy table, convert the
# existing Delta [%] formulas (C4:C12) into a shared formula group, and
# move the active selection to B14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row to the "Tableau1" table (extends it to A3:C13) ---
$lo = $ws.ListObjects.Item("Tableau1")
$newRow = $lo.ListRows.Add()
$newRow.Range.Item(1).Value = 150000
$newRow.Range.Item(2).Value = 150093.81
$newRow.Range.Item(2).NumberFormat = "0.000"
$newRow.Range.Item(3).Formula = "=B13/(A13/100)-100"
$newRow.Range.Item(3).NumberFormat = "0.000"

# --- Re-enter the Delta [%] formula across the original rows so Excel
#     stores it as a single shared formula (C4:C12) like it does when a
#     formula is typed once and filled down a column ---
$ws.Range("C4:C12").Formula = "=B4/(A4/100)-100"

# --- Update the selected cell on the sheet ---
$ws.Range("B14").Select()
